$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-04-17 Wednesday"; new="2024-04-18 Thursday"},
    @{old="49×38=1862"; new="67×35=2345"},
    @{old="59×12=708"; new="25×85=2125"},
    @{old="48×47=2256"; new="35×52=1820"},
    @{old="62×78=4836"; new="50×49=2450"},
    @{old="62×49=3038"; new="83×27=2241"},
    @{old="68×33=2244"; new="53×21=1113"},
    @{old="87×47=4089"; new="73×52=3796"},
    @{old="67×94=6298"; new="87×86=7482"},
    @{old="62×51=3162"; new="26×42=1092"},
    @{old="84×11=924"; new="71×91=6461"},
    @{old="17×21=357"; new="34×99=3366"},
    @{old="73×79=5767"; new="50×23=1150"},
    @{old="42×18=756"; new="37×28=1036"},
    @{old="27×62=1674"; new="28×94=2632"},
    @{old="90×49=4410"; new="64×77=4928"},
    @{old="47×72=3384"; new="29×98=2842"},
    @{old="57×79=4503"; new="66×93=6138"},
    @{old="61×16=976"; new="88×82=7216"},
    @{old="87×24=2088"; new="76×37=2812"},
    @{old="24×15=360"; new="46×99=4554"},
    @{old="11×64=704"; new="85×19=1615"},
    @{old="11×27=297"; new="28×81=2268"},
    @{old="37×70=2590"; new="40×46=1840"},
    @{old="80×16=1280"; new="95×45=4275"},
    @{old="77×16=1232"; new="20×48=960"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
